$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The sheet is a sorted price list. Two new product rows need to be inserted
# at their correct alphabetical position:
#   - "GARAMYCIN 0.1% OINT. 15 GM"      -> becomes new row 17
#   - "PANTHENOL 2% TOPICAL CREAM 50 GM" -> becomes new row 29 (after the
#     first insert shifts everything below row 17 down by one)
# Every row below each insertion point shifts down by one; the "#" (A column)
# sequence, the totals row and the footer row all move down accordingly.
# ---------------------------------------------------------------------------

# Insert the GARAMYCIN row at row 17, cloning the formatting of row 16
# (the row immediately above it) so styles/borders/number-formats match.
$ws.Rows("17:17").Insert()
$ws.Range("A16:N16").Copy($ws.Range("A17:N17"))
$ws.Rows("17:17").RowHeight = 25.5

# Insert the PANTHENOL 50 GM row at row 30, right after the existing
# "PANTHENOL 2% TOPICAL CREAM 20 GM" row which is now at row 29 (it was row
# 28 before the first insert shifted it down by one).
$ws.Rows("30:30").Insert()
$ws.Range("A29:N29").Copy($ws.Range("A30:N30"))
$ws.Rows("30:30").RowHeight = 24.75

# ---------------------------------------------------------------------------
# Rewrite the full data block (rows 4-40) with its final values. Column A is
# just the running "#" counter, column B/H hold the (merged) product name and
# usage-ratio text, L is the balance qty and N is the sale-price ratio.
# ---------------------------------------------------------------------------
$rows = @(
  @(4,  1,  'ALPHINTERN 30 F.C.TABS', '4:2', 29, 0.33),
  @(5,  2,  'ALVEOLIN-P SYRUP 100 ML', '1:0', 50, 1),
  @(6,  3,  'ANTINAL 220MG/5ML 60ML SUSP.', '1:0', 48, 2),
  @(7,  4,  'ANTODINE20    6 AMP', '0:2', 13, 0.17),
  @(8,  5,  'APEXIDONE 4MG 30 F.C.TAB.', '0:2', 42, 0.33),
  @(9,  6,  'ATROVENT 250MCG/2ML 20 UNIT DOSE VIAL', '1:19', 42.9, 0.15),
  @(10, 7,  'BECOZYME I.M./I.V. 12 AMP', '0:7', 10, 0.08),
  @(11, 8,  'CERVITAM 20 CAPS.', '0:1', 59, 0.5),
  @(12, 9,  'DEPAKINE CHRONO 500MG 30 SCORED PROLONGED REL. F.C. TAB.', '1:0', 139.68, 1),
  @(13, 10, 'DEXAMETHASONE INAD PHARMA 3 AMP', '6:0', 12, 0.33),
  @(14, 11, 'DIAMICRON MR 30 MG 30 TAB.', '1:1', 26.67, 0.33),
  @(15, 12, 'DOLIPRANE 1 GM 15 TABS.', '10:1', 32, 0.67),
  @(16, 13, 'DRAMENEX 50MG 20 TABS.', '1:0', 28, 1),
  @(17, 14, 'GARAMYCIN 0.1% OINT. 15 GM', '1:0', 44, 2),
  @(18, 15, 'GOURYST 0.5 MG 100 TABS.', '0:6', 19, 0.1),
  @(19, 16, 'KETOLAC 30MG/2ML 5 AMP. FOR I.M./I.V. INF.', '2:0', 12, 0.2),
  @(20, 17, 'MEGAFEN-N 100MG/5ML SUSP. 120 ML', '1:0', 35, 1),
  @(21, 18, 'MIXDERM CREAM 30 GM', '2:0', 41, 1),
  @(22, 19, 'MUCOPHYLLINE SYRUP 125 ML', '10:0', 50, 1),
  @(23, 20, 'MUCOSTA 100MG 20 TAB', '3:1', 57, 0.5),
  @(24, 21, 'NEUROVIT 6 I.M. AMPS', '3:4', 11, 0.17),
  @(25, 22, 'NEXIUM 20MG 28 F.C. TAB.', '0:0', 332, 1),
  @(26, 23, 'OPLEX-N SYRUP 125ML', '5:0', 31, 1),
  @(27, 24, 'PANADOL ADVANCE 500 MG 48 TABLETS', '2:3', 23, 0.25),
  @(28, 25, 'PANTHENOL 2% TOPICAL CREAM 20 GM', '4:0', 32, 1),
  @(29, 26, 'PANTHENOL 2% TOPICAL CREAM 50 GM', '1:0', 80, 1),
  @(30, 27, 'PULMICORT 0.25MG/ML 20 NEBULIZER VIAL SUSP.', '0:19', 169.2, 0.3),
  @(31, 28, 'SPASMO-DIGESTIN 30 TABS.', '3:1', 73.32, 1),
  @(32, 29, 'SUGARLO PLUS 50/1000MG 30 F.C. TABS', '2:1', 45.5, 0.33),
  @(33, 30, 'TELFAST 120MG 20 F.C. TAB', '0:0', 116, 1),
  @(34, 31, 'TUSSISTOP 60 MG 20 TABS.', '0:0', 30, 0.5),
  @(35, 32, 'VOLTAREN 75MG/3ML 3 AMP.', '4:1', 17, 0.33),
  @(36, 33, 'YEAST MEPACO 60 TABS', '1:0', 60, 1),
  @(37, 34, 'جهاز محلول ', '3:0', 20, 1),
  @(38, 35, 'سرنجات 3 سم', '-2:0', 4, 2),
  @(39, 36, 'سرنجات 5 سم', '-1:0', 2, 1),
  @(40, 37, 'محلول خليط', '3:0', 27, 1)
)

foreach ($row in $rows) {
  $r = $row[0]
  $ws.Range("A$r").Value = $row[1]
  $ws.Range("B$r").Value = $row[2]
  $ws.Range("H$r").Value = $row[3]
  $ws.Range("L$r").Value = $row[4]
  $ws.Range("N$r").Value = $row[5]
}

# Totals row (was row 39, now row 41): grand total balance-qty sum.
$ws.Range("K41").Value = 1863.27

Write-Output "edit complete"
